# Regenerate save_data: replace the "K" column (G) values with the
# recomputed strikeout counts from the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 2
    29 = 2
    30 = 1
    31 = 2
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    44 = 3
    45 = 2
    46 = 1
    47 = 0
    48 = 0
    50 = 1
    51 = 0
    52 = 2
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 1
    60 = 0
    61 = 1
    62 = 0
    63 = 1
    65 = 1
    66 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
